$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Formula = "'27.198.49"
$ws.Cells.Item(2, 5).Value = "  +0.05%  "
$ws.Cells.Item(3, 4).Formula = "'1.631.76"
$ws.Cells.Item(3, 5).Value = "  -0.82%  "
$ws.Cells.Item(4, 5).Value = "  +0.04%  "
$ws.Cells.Item(5, 4).Formula = "'216.62"
$ws.Cells.Item(5, 5).Value = "  -0.49%  "
$ws.Cells.Item(6, 4).Formula = "'0.517"
$ws.Cells.Item(6, 5).Value = "  +1.38%  "
$ws.Cells.Item(7, 5).Value = "  +0.03%  "
$ws.Cells.Item(8, 5).Value = "  -0.52%  "
$ws.Cells.Item(9, 4).Formula = "'0.0625"
$ws.Cells.Item(9, 5).Value = "  -0.88%  "
$ws.Cells.Item(10, 4).Formula = "'20.31"
$ws.Cells.Item(10, 5).Value = "  +1.52%  "
$ws.Cells.Item(11, 4).Formula = "'0.0848"
$ws.Cells.Item(11, 5).Value = "  +0.66%  "
$ws.Cells.Item(12, 4).Formula = "'1.632.45"
$ws.Cells.Item(12, 5).Value = "  -1.27%  "
$ws.Cells.Item(13, 4).Formula = "'4.12"
$ws.Cells.Item(13, 5).Value = "  -0.01%  "
$ws.Cells.Item(14, 4).Formula = "'0.543"
$ws.Cells.Item(14, 5).Value = "  +0.83%  "
$ws.Cells.Item(15, 4).Formula = "'27.147.10"
$ws.Cells.Item(15, 5).Value = "  +0.00%  "
$ws.Cells.Item(16, 4).Formula = "'64.70"
$ws.Cells.Item(16, 5).Value = "  -4.02%  "
$ws.Cells.Item(17, 5).Value = "  -0.74%  "
$ws.Cells.Item(18, 4).Formula = "'215.09"
$ws.Cells.Item(18, 5).Value = "  -1.74%  "
$ws.Cells.Item(19, 5).Value = "  +0.16%  "
$ws.Cells.Item(20, 5).Value = "  +0.74%  "
$ws.Cells.Item(21, 5).Value = "  -0.82%  "
$ws.Cells.Item(22, 4).Formula = "'2.49"
$ws.Cells.Item(22, 5).Value = "  -0.19%  "
$ws.Cells.Item(23, 5).Value = "  -1.03%  "
$ws.Cells.Item(24, 4).Formula = "'148.35"
$ws.Cells.Item(24, 5).Value = "  +0.52%  "
$ws.Cells.Item(25, 5).Value = "  +0.05%  "
$ws.Cells.Item(26, 4).Formula = "'7.29"
$ws.Cells.Item(26, 5).Value = "  -1.77%  "
$ws.Cells.Item(27, 5).Value = "  +0.01%  "
$ws.Cells.Item(28, 5).Value = "  -1.10%  "
$ws.Cells.Item(29, 4).Formula = "'0.0505"
$ws.Cells.Item(29, 5).Value = "  -0.10%  "
$ws.Cells.Item(30, 5).Value = "  -0.74%  "
$ws.Cells.Item(31, 5).Value = "  +0.62%  "
$ws.Cells.Item(32, 5).Value = "  -0.50%  "
$ws.Cells.Item(33, 4).Formula = "'1.313.09"
$ws.Cells.Item(33, 5).Value = "  +3.73%  "
$ws.Cells.Item(34, 5).Value = "  -1.53%  "
$ws.Cells.Item(35, 4).Formula = "'2.45"
$ws.Cells.Item(35, 5).Value = "  +0.15%  "
$ws.Cells.Item(36, 5).Value = "  -1.53%  "
$ws.Cells.Item(37, 4).Formula = "'0.850"
$ws.Cells.Item(37, 5).Value = "  +1.13%  "
$ws.Cells.Item(38, 4).Formula = "'0.540"
$ws.Cells.Item(39, 5).Value = "  +0.07%  "
$ws.Cells.Item(40, 4).Formula = "'2.26"
$ws.Cells.Item(40, 5).Value = "  +1.62%  "
$ws.Cells.Item(41, 4).Formula = "'0.806"
$ws.Cells.Item(41, 5).Value = "  -0.53%  "
$ws.Cells.Item(42, 4).Formula = "'63.67"
$ws.Cells.Item(42, 5).Value = "  +2.28%  "
$ws.Cells.Item(43, 4).Formula = "'1.768.93"
$ws.Cells.Item(43, 5).Value = "  -0.93%  "
$ws.Cells.Item(44, 5).Value = "  -2.80%  "
$ws.Cells.Item(45, 4).Formula = "'90.73"
$ws.Cells.Item(45, 5).Value = "  -1.15%  "
$ws.Cells.Item(46, 4).Formula = "'1.60"
$ws.Cells.Item(46, 5).Value = "  +0.00%  "
$ws.Cells.Item(47, 2).Value = "WEMIXToken"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(47, 4).Formula = "'0.799"
$ws.Cells.Item(47, 5).Value = "  +18.26%  "
$ws.Cells.Item(48, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(48, 4).Formula = "'0.0₆0102"
$ws.Cells.Item(48, 5).Value = "  -5.21%  "
$ws.Cells.Item(49, 4).Formula = "'0.0517"
$ws.Cells.Item(49, 5).Value = "  +1.00%  "
$ws.Cells.Item(50, 4).Formula = "'7.52"
$ws.Cells.Item(50, 5).Value = "  -2.34%  "
$ws.Cells.Item(51, 5).Value = "  -0.16%  "
